$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'256.42"
$ws.Range("E2").Value = "'-0.42%"
$ws.Range("G2").Value = "'8"

# Row 3
$ws.Range("D3").Value = "'27.10"
$ws.Range("E3").Value = "'-0.45%"
$ws.Range("G3").Value = "'8"

# Row 4
$ws.Range("D4").Value = "'4.493"
$ws.Range("E4").Value = "'-5.89%"
$ws.Range("G4").Value = "'8"

# Row 5
$ws.Range("D5").Value = "'0.05890"
$ws.Range("E5").Value = "'-1.20%"
$ws.Range("G5").Value = "'8"

# Row 6
$ws.Range("D6").Value = "'6.610"
$ws.Range("E6").Value = "'-0.83%"
$ws.Range("G6").Value = "'8"

# Row 7
$ws.Range("D7").Value = "'0.8501"
$ws.Range("E7").Value = "'-2.27%"
$ws.Range("G7").Value = "'8"

# Row 8
$ws.Range("D8").Value = "'0.9227"
$ws.Range("E8").Value = "'-3.97%"
$ws.Range("G8").Value = "'8"

# Row 9
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.01034"
$ws.Range("E9").Value = "'1,596.24%"
$ws.Range("G9").Value = "'8"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1376"
$ws.Range("E10").Value = "'-2.20%"
$ws.Range("G10").Value = "'8"

# Row 11
$ws.Range("B11").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C11").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D11").Value = "'0.04314"
$ws.Range("E11").Value = "'16.26%"
$ws.Range("G11").Value = "'8"

# Row 12
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07008"
$ws.Range("E12").Value = "'-2.06%"
$ws.Range("G12").Value = "'8"

# Row 13
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.03083"
$ws.Range("E13").Value = "'-2.85%"
$ws.Range("G13").Value = "'8"

# Row 14
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.09095"
$ws.Range("E14").Value = "'-1.70%"
$ws.Range("G14").Value = "'8"

# Row 15
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001532"
$ws.Range("E15").Value = "'-0.48%"
$ws.Range("G15").Value = "'8"

# Row 16
$ws.Range("D16").Value = "'0.006042"
$ws.Range("E16").Value = "'-0.86%"
$ws.Range("G16").Value = "'8"

# Row 17
$ws.Range("D17").Value = "'3.470"
$ws.Range("E17").Value = "'-0.31%"
$ws.Range("G17").Value = "'8"

# Row 18
$ws.Range("D18").Value = "'3.170"
$ws.Range("E18").Value = "'-0.49%"
$ws.Range("G18").Value = "'8"

# Row 19
$ws.Range("D19").Value = "'2.194"
$ws.Range("E19").Value = "'-1.12%"
$ws.Range("G19").Value = "'8"

# Row 20
$ws.Range("D20").Value = "'0.3026"
$ws.Range("E20").Value = "'-3.37%"
$ws.Range("G20").Value = "'8"

# Row 21
$ws.Range("E21").Value = "'-1.56%"
$ws.Range("G21").Value = "'8"

# Row 22
$ws.Range("D22").Value = "'3.899"
$ws.Range("E22").Value = "'2.31%"
$ws.Range("G22").Value = "'8"

# Row 23
$ws.Range("D23").Value = "'0.04270"
$ws.Range("E23").Value = "'1.10%"
$ws.Range("G23").Value = "'8"

# Row 24
$ws.Range("D24").Value = "'0.001219"
$ws.Range("E24").Value = "'-0.12%"
$ws.Range("G24").Value = "'8"

# Row 25
$ws.Range("D25").Value = "'0.004301"
$ws.Range("E25").Value = "'-4.40%"
$ws.Range("G25").Value = "'8"

# Row 26
$ws.Range("D26").Value = "'0.0001199"
$ws.Range("E26").Value = "'-0.02%"
$ws.Range("G26").Value = "'8"

# Row 27
$ws.Range("G27").Value = "'8"

# Row 28
$ws.Range("G28").Value = "'8"

# Row 29
$ws.Range("G29").Value = "'8"

# Row 30
$ws.Range("G30").Value = "'8"

# Row 31
$ws.Range("G31").Value = "'8"

# Row 32
$ws.Range("G32").Value = "'8"

# Row 33
$ws.Range("G33").Value = "'8"

# Row 34
$ws.Range("G34").Value = "'8"

# Row 35
$ws.Range("G35").Value = "'8"

# Row 36
$ws.Range("G36").Value = "'8"

# Row 37
$ws.Range("G37").Value = "'8"

# Row 38
$ws.Range("G38").Value = "'8"

# Row 39
$ws.Range("G39").Value = "'8"

# Row 40
$ws.Range("D40").Value = "'0.03807"
$ws.Range("E40").Value = "'-0.36%"
$ws.Range("G40").Value = "'8"

# Row 41
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.006257"
$ws.Range("E41").Value = "'0.14%"
$ws.Range("G41").Value = "'8"

# Row 42
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1099"
$ws.Range("E42").Value = "'0.01%"
$ws.Range("G42").Value = "'8"

# Row 43
$ws.Range("D43").Value = "'0.002199"
$ws.Range("E43").Value = "'-2.33%"
$ws.Range("G43").Value = "'8"

# Row 44
$ws.Range("D44").Value = "'0.01399"
$ws.Range("E44").Value = "'31.96%"
$ws.Range("G44").Value = "'8"

# Row 45
$ws.Range("D45").Value = "'0.00005351"
$ws.Range("E45").Value = "'-2.69%"
$ws.Range("G45").Value = "'8"

# Row 46
$ws.Range("E46").Value = "'-0.02%"
$ws.Range("G46").Value = "'8"

# Row 47
$ws.Range("D47").Value = "'0.05377"
$ws.Range("G47").Value = "'8"

# Row 48
$ws.Range("E48").Value = "'10,589.89%"
$ws.Range("G48").Value = "'8"

# Row 49
$ws.Range("D49").Value = "'0.00002099"
$ws.Range("E49").Value = "'-0.02%"
$ws.Range("G49").Value = "'8"

# Row 50
$ws.Range("D50").Value = "'0.0001999"
$ws.Range("E50").Value = "'-0.02%"
$ws.Range("G50").Value = "'8"

# Row 51
$ws.Range("G51").Value = "'8"
